$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") values for rows 2-24 should move from serial date
# 46060 (2026-02-07) to 46061 (2026-02-08).
for ($row = 2; $row -le 24; $row++) {
    $ws.Cells.Item($row, 3).Value = Get-Date -Year 2026 -Month 2 -Day 8 -Hour 0 -Minute 0 -Second 0
}
